$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 87, shifting existing rows 87.. down by one
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with data
$ws.Cells.Item(87, 1).Value = 4
$ws.Cells.Item(87, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(87, 3).Value = "Los Lagos"
$ws.Cells.Item(87, 4).Value = 44512
$ws.Cells.Item(87, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(87, 5).Value = 10
$ws.Cells.Item(87, 6).Value = 100112003
$ws.Cells.Item(87, 7).Value = "Ajo"
$ws.Cells.Item(87, 8).Value = "Chino"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 240
$ws.Cells.Item(87, 11).Value = 21000
$ws.Cells.Item(87, 12).Value = 22000
$ws.Cells.Item(87, 13).Value = 21500
$ws.Cells.Item(87, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(87, 15).Value = "China"
$ws.Cells.Item(87, 16).Value = 2150
$ws.Cells.Item(87, 17).Value = 10
$ws.Cells.Item(87, 18).Value = "Hortaliza"
